$d = $word.ActiveDocument

$pairs = @(
    @("30×75=", "71×67="),
    @("84×41=", "72×12="),
    @("86×82=", "88×61="),
    @("87×95=", "81×49="),
    @("48×54=", "80×75="),
    @("45×64=", "62×15="),
    @("55×97=", "19×67="),
    @("76×33=", "77×81="),
    @("75×13=", "99×44="),
    @("52×15=", "87×56="),
    @("65×84=", "60×93="),
    @("33×64=", "24×21="),
    @("36×25=", "51×46="),
    @("82×60=", "79×71="),
    @("17×91=", "21×19="),
    @("15×14=", "55×27="),
    @("97×55=", "53×35="),
    @("28×94=", "38×48="),
    @("62×52=", "99×66="),
    @("36×45=", "69×90="),
    @("15×79=", "97×74="),
    @("84×96=", "91×76="),
    @("69×88=", "80×39="),
    @("16×46=", "67×11="),
    @("75×41=", "40×16=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
